$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (by sheet row number) that get updated "Number_of_Inclusions" (B)
# and recomputed "Number_of_Inclusions_per_Nucleus" (D = B / C), now
# using the 3rd quartile instead of the mean.
$updates = @{
    4  = 21
    5  = 2
    6  = 84
    7  = 13
    8  = 29
    11 = 0
    12 = 104
    13 = 0
    19 = 0
    22 = 5
    23 = 32
}

foreach ($row in $updates.Keys) {
    $newB = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $newB

    $nuclei = $ws.Cells.Item($row, 3).Value()
    if ($nuclei -ne 0) {
        $ws.Cells.Item($row, 4).Value = $newB / $nuclei
    } else {
        $ws.Cells.Item($row, 4).Value = 0
    }
}
